# atualizei dados da add
# Updates the June-2025 daily revenue figures and inserts the missing
# 17/06/2025 row (previously absent from the sheet), shifting every
# subsequent row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revised totals for existing June rows (day 9, 11, 13, 16)
$ws.Range("B7").Value  = 10652.15
$ws.Range("B9").Value  = 19810.51
$ws.Range("B11").Value = 9230.379999999999
$ws.Range("B12").Value = 24791.07

# Insert a new row right after the current row 12 (day 16) to hold the
# previously-missing day-17 record; this shifts rows 13:72 down to 14:73.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row with the day-17 June data.
$ws.Range("A13").Value = 17
$ws.Range("B13").Value = 13628.15
$ws.Range("C13").Value = 6
$ws.Range("D13").Value = 2025
$ws.Range("E13").Value = "06/2025"
